$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DM_Stat (C) and P_Value (D) columns with corrected values
$ws.Range("C2").Value = 0.9986931695954243
$ws.Range("D2").Value = 0.3249983597280082

$ws.Range("C3").Value = 0.9020191021013145
$ws.Range("D3").Value = 0.3733926597242712

$ws.Range("C4").Value = -0.008218464787892052
$ws.Range("D4").Value = 0.9934907216136308

$ws.Range("C5").Value = -1.619712381982112
$ws.Range("D5").Value = 0.1145339259482292

$ws.Range("C6").Value = -0.6012970551071446
$ws.Range("D6").Value = 0.5516304305335686

$ws.Range("C7").Value = -1.431138537325558
$ws.Range("D7").Value = 0.1615185792220515

$ws.Range("C8").Value = -1.64200636485271
$ws.Range("D8").Value = 0.1098075962741822
$ws.Range("G8").Value = "No"

$ws.Range("C9").Value = -1.189886919214864
$ws.Range("D9").Value = 0.2423335535766857

$ws.Range("C10").Value = -1.604068123660045
$ws.Range("D10").Value = 0.1179493680421706
$ws.Range("G10").Value = "No"

$ws.Range("C11").Value = -1.265508176158052
$ws.Range("D11").Value = 0.2142931672162649

$wb.Save()
